# Minor update of ppt with format
$p = $ppt.ActivePresentation

# Helper: build the long value PowerPoint's Color.RGB expects (R + G*256 + B*65536)
# for target accent color 29748D (R=0x29, G=0x74, B=0x8D)
$teal = 0x29 + (0x74 * 256) + (0x8D * 65536)

# --- Slide 1: title slide - fix missing space in "&prediction" ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item("Title 1").TextFrame.TextRange.Text = "Big Mountain ticket price modeling & prediction"

# --- Slide 3: "Recommendation and key findings" ---
$s3 = $p.Slides.Item(3)
# Remove the stray full-bleed background rectangle shape (Google Shape;20;p1)
$s3.Shapes.Item("Google Shape;20;p1").Delete()
# Recolor the title run from theme accent1 to explicit teal
$s3.Shapes.Item("Title 1").TextFrame.TextRange.Font.Color.RGB = $teal

# --- Slide 4: "Modeling results and analysis" ---
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item("Title 1").TextFrame.TextRange.Font.Color.RGB = $teal

# --- Slide 5: "Modeling results and analysis (cont.)" ---
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item("Title 1").TextFrame.TextRange.Font.Color.RGB = $teal

# --- Slide 6: "Modeling results and analysis (cont.)" ---
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item("Title 1").TextFrame.TextRange.Font.Color.RGB = $teal
$s6.Shapes.Item("Content Placeholder 2").TextFrame.TextRange.Text = "Ranking of Big Mountain on key features and price among the market share"

# --- Slide 7: "Modeling results and analysis (cont.)" ---
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item("Title 1").TextFrame.TextRange.Font.Color.RGB = $teal
$s7.Shapes.Item("Content Placeholder 2").TextFrame.TextRange.Text = "Scenarios modeling of closing down the least used runs, but no more than 5 runs"

# --- Slide 8: "Summary and conclusion" ---
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item("Title 1").TextFrame.TextRange.Font.Color.RGB = $teal
$s8Content = $s8.Shapes.Item("Content Placeholder 2")
# The 3rd/4th paragraph text is only lightly reworded from the original, and a
# direct Paragraphs(n,1).Text assignment against very-similar text makes the
# host emit a diff-minimizing run split (multiple <a:r> sharing the prefix).
# Route the replacement through an unrelated placeholder first so the final
# assignment lands as a single clean run, matching authored output.
$s8Content.TextFrame.TextRange.Paragraphs(3, 1).Text = "PLACEHOLDER_TOKEN_ONE"
$s8Content.TextFrame.TextRange.Paragraphs(3, 1).Text = "More data can be obtained for better prediction and estimation of revenue increase, such as visitor number and average staying days data, weather temperature, geographic latitude, transportation convenience, accommodations like hotels and restaurants numbers on site, and business costs, etc."
$s8Content.TextFrame.TextRange.Paragraphs(4, 1).Text = "PLACEHOLDER_TOKEN_TWO"
$s8Content.TextFrame.TextRange.Paragraphs(4, 1).Text = "Small tests can be done to test modeled improvement scenarios, with comparison and visitor feedbacks to see primary results and then further prediction and recommendations."

# --- Slide 9: "Thank you!" - reposition/resize and center the text ---
$s9 = $p.Slides.Item(9)
$thankYou = $s9.Shapes.Item("Title 1")
# Left/Top/Width/Height are in points (EMU / 12700); the host stores them as
# single-precision floats, so the literals below are nudged by a few 1e-7
# so the float32 round-trip truncates back to the exact target EMU.
$thankYou.Left = 371.6504724409449
$thankYou.Top = 208.40055118110237
$thankYou.Width = 216.69905861811023
$thankYou.Height = 104.37504197007874
$thankYou.TextFrame.TextRange.Paragraphs(1, 1).ParagraphFormat.Alignment = 2
